$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D5").Value = 1.035825651653654
$ws.Range("E5").Value = 18.2677832683071
$ws.Range("F5").Value = 62.67199190538536
$ws.Range("G5").Value = 19.06022482630755
$ws.Range("H5").Value = 80.93977517369245
$ws.Range("I5").Value = 19.06022482630755
$ws.Range("J5").Value = 0.3869388155426705
$ws.Range("K5").Value = 6.655100096637757
$ws.Range("L5").Value = 86.09483034719469
$ws.Range("M5").Value = 7.250069556167555
$ws.Range("N5").Value = 92.74993044383244
$ws.Range("O5").Value = 7.250069556167555

$ws.Range("C6").Value = 8.03286593
$ws.Range("J6").Value = 0.3175744464262313
$ws.Range("K6").Value = 2.283018344493816
$ws.Range("L6").Value = 67.55088713438603
$ws.Range("M6").Value = 30.16609452112015
$ws.Range("N6").Value = 69.83390547887984
$ws.Range("O6").Value = 30.16609452112015

$ws.Range("D22").Value = 0.05364998000000001
$ws.Range("E22").Value = 49.92385831271513
$ws.Range("F22").Value = 38.86139379735091
$ws.Range("G22").Value = 11.21474788993398
$ws.Range("H22").Value = 88.78525211006604
$ws.Range("I22").Value = 11.21474788993398

$ws.Range("C24").Value = 80.28050343000002
$ws.Range("D24").Value = 46.17195226549329
$ws.Range("E24").Value = 29.05515303058552
$ws.Range("F24").Value = 50.87049853954008
$ws.Range("G24").Value = 20.07434842987439
$ws.Range("H24").Value = 79.92565157012559
$ws.Range("I24").Value = 20.07434842987439
$ws.Range("J24").Value = 15.45779683901184
$ws.Range("K24").Value = 20.18311283805593
$ws.Range("L24").Value = 48.98125251884212
$ws.Range("M24").Value = 30.83563464310195
$ws.Range("N24").Value = 69.16436535689806
$ws.Range("O24").Value = 30.83563464310195
